# Scene.xlsx: XML schema for the Scene ini table changed -- the SceneID
# column was dropped, a FilePath column was promoted to the front, and the
# ID values switched from "SceneN" labels to plain numeric-looking index
# strings ("0","1","2"). Column order is now alphabetical:
#   FilePath, ID, MaxGroup, MaxGroupPlayers, RelivePos, SceneName,
#   SceneShowName, SoundList, Width
# Also: the source XML connection moved from Excel_Ini\Scene.xml to
# Ini\NPC\Scene.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Start from a clean slate for the A1:J4 block so no stale number formats
# leak into the rebuilt layout.
$ws.Range("A1:J4").ClearFormats()

# ---- header row -----------------------------------------------------
$headers = @("FilePath", "ID", "MaxGroup", "MaxGroupPlayers", "RelivePos", "SceneName", "SceneShowName", "SoundList", "Width")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

# ---- data rows --------------------------------------------------------
# FilePath, ID, MaxGroup, MaxGroupPlayers, RelivePos, SceneName, SceneShowName, SoundList, Width
$data = @(
    @("../../NFDataCfg/Ini/NFZoneServer/Scene/PioneerNoob/", "0", 10, 200, "-0.7,1.5,-1.64", "Login", "NewerScene", $null, 500),
    @("../../NFDataCfg/Ini/NFZoneServer/Scene/PioneerNoob/", "1", 10, 200, "-0.7,1.5,-1.64", "Stage001", "NewerScene", $null, 500),
    @("../../NFDataCfg/Ini/NFZoneServer/Scene/RebellerNoob/", "2", 10, 200, "-0.7,1.5,-1.64", "Stage001", "NewerScene", $null, 500)
)

$textCols = @(1, 2, 5, 6, 7, 8)   # FilePath, ID, RelivePos, SceneName, SceneShowName, SoundList -> text format
$numCols = @(3, 4, 9)             # MaxGroup, MaxGroupPlayers, Width -> general/number format

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowIdx = $r + 2
    $rowVals = $data[$r]
    foreach ($c in $textCols) {
        $cell = $ws.Cells.Item($rowIdx, $c)
        $cell.NumberFormat = "@"
        $val = $rowVals[$c - 1]
        if ($null -ne $val) {
            $cell.Value = $val
        }
    }
    foreach ($c in $numCols) {
        $cell = $ws.Cells.Item($rowIdx, $c)
        $cell.Value = $rowVals[$c - 1]
    }
}

# ---- drop the now-unused column J (old SoundList slot) ---------------
$ws.Columns.Item(10).Delete()

# ---- resync the XML-mapped table to the new A1:I4 extent -------------
$lo.Resize($ws.Range("A1:I4"))

# ---- column widths (character units) ----------------------------------
$widths = @{1 = 58.25; 2 = 7.5; 3 = 12.75; 4 = 21.5; 5 = 16.125; 6 = 14; 7 = 19; 8 = 14; 9 = 9.375}
foreach ($c in $widths.Keys) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c]
}

# ---- selection / dimension ---------------------------------------------
$ws.Range("F1:F1048576").Select()

# ---- source XML connection path (best effort) --------------------------
$qt = $lo.QueryTable
if ($qt) {
    $qt.SourceDataFile = "D:\NoahGameFrame\trunk\_Out\Server\NFDataCfg\Ini\NPC\Scene.xml"
}

Write-Host "Scene.xlsx restructured: FilePath/ID/MaxGroup/MaxGroupPlayers/RelivePos/SceneName/SceneShowName/SoundList/Width"
